# Reorder the "Recorded By" (column G) values so that any "System"/"system"
# token is moved to the front of the comma-separated list, while the
# relative order of the remaining tokens is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -gt 0) {
        $newOrder = $systemParts + $otherParts
        $newVal = [string]::Join(", ", $newOrder)
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
